# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# Sheet "OFF" (first sheet) - row 3 ("R" row) gets updated Short/Deep Att/Comp/Int counts
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 268
$wsOff.Range("C3").Value = 179
$wsOff.Range("D3").Value = 90
$wsOff.Range("E3").Value = 40

# Sheet "DEF" (second sheet) - row 3 ("R" row) gets updated Short/Deep Att/Comp/Int counts
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 499
$wsDef.Range("C3").Value = 372
$wsDef.Range("D3").Value = 79
$wsDef.Range("E3").Value = 38
$wsDef.Range("G3").Value = 3
